# Fill in the Date / Opponent-guild / Attack-or-defense columns (H, I, J)
# for the newly-added match rows 267-314, mirroring the values already
# present on the preceding rows (240-266) for the same guild-war date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$fillRange = $ws.Range("H267:J314")
$fillRange.HorizontalAlignment = -4108   # xlCenter, matches the rest of the sheet's "s=1" style
$fillRange.VerticalAlignment = -4108     # xlCenter

for ($r = 267; $r -le 314; $r++) {
    $ws.Cells.Item($r, 8).Value = 260107        # H: 날짜 (date)
    $ws.Cells.Item($r, 9).Value = "푸른달"        # I: 상대 길드 (opponent guild)
    $ws.Cells.Item($r, 10).Value = "공격"         # J: 기준 (attack/defense)
}

# Restore the recorded scroll position / active selection for this sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 266
$win.ScrollColumn = 1
$ws.Range("N310").Select()
